$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(21.52698182277445, 3.023968074486533, 3.904400907923463, 10.77759118223927, 69.53616316249635, 0, 7.344005520526261, 0, 11.27624633093289, 17.68412328599214, 11.10012192625136),
    @(21.54617254014096, 2.974724868052053, 3.780204195458183, 10.81090140214564, 68.71104351508677, 0, 7.344005520526261, 0, 11.26839204990293, 17.71754574697365, 11.14497077895294),
    @(21.56471972667961, 2.946699213255214, 3.700840707651576, 10.83268889960513, 68.20216886422506, 0, 7.344005520526261, 0, 11.26399002304225, 17.74346305139048, 11.17448627261451),
    @(21.57397614888358, 2.935857803289564, 3.667736998744604, 10.84190398344511, 67.9943717979898, 0, 7.344005520526261, 0, 11.26230208162477, 17.75537804697915, 11.18701203763196),
    @(21.57561564859246, 2.934093227693434, 3.662194687350429, 10.84345448954115, 67.95984553633438, 0, 7.344005520526261, 0, 11.26202820225961, 17.75743816180599, 11.18912202337118),
    @(21.5648376899278, 2.94655062724572, 3.700397320412901, 10.83281181386313, 68.19936798461718, 0, 7.344005520526261, 0, 11.26396682978006, 17.74361826586571, 11.17465318252005),
    @(21.53219366327373, 3.006545158122624, 3.862233136271364, 10.78879999072991, 69.25221172474714, 0, 7.344005520526261, 0, 11.27345099982036, 17.69452561918748, 11.11517557995334),
    @(21.52191954814764, 3.2038286111617, 4.154229233018861, 10.71304731982904, 71.2929266933364, 0, 7.344005520526261, 0, 11.29537928146186, 17.64119333070693, 11.0142114707339),
    @(21.54716606557558, 3.445411907778826, 4.352469972263616, 10.66377474212559, 72.76961114860046, 0, 7.344005520526261, 0, 11.31351182643871, 17.62832200886975, 10.94955415014473),
    @(21.56575678325919, 3.549603523146622, 4.438977206417737, 10.64273455817717, 73.43480147946703, 0, 7.344005520526261, 0, 11.32219778337085, 17.62819591828743, 10.92220065632671),
    @(21.57381518127434, 3.588237162318769, 4.471197135727511, 10.63496396522376, 73.68561018233571, 0, 7.344005520526261, 0, 11.32554956237788, 17.62897218307044, 10.91213828440935),
    @(21.57203442844455, 3.579953257183808, 4.464282095176294, 10.63662875728736, 73.63164448069628, 0, 7.344005520526261, 0, 11.32482491782295, 17.62876835703995, 10.91429224562158),
    @(21.56639936350234, 3.552798401456598, 4.441638805175368, 10.64209132546392, 73.45545819194629, 0, 7.344005520526261, 0, 11.32247228352154, 17.62824326911458, 10.92136689251632),
    @(21.56308024331099, 3.53605827679561, 4.427698729417267, 10.64546292259101, 73.34739364873704, 0, 7.344005520526261, 0, 11.32103936966248, 17.62802894155954, 10.92573883033895),
    @(21.54609385001092, 3.43848781927904, 4.346741694147365, 10.66517735624739, 72.72599524839266, 0, 7.344005520526261, 0, 11.31295295480438, 17.62844556435533, 10.9513831768762),
    @(21.53749124913137, 3.377169094787589, 4.296128578311471, 10.67762295251331, 72.34300574043628, 0, 7.344005520526261, 0, 11.30810401767048, 17.63016885814851, 10.96764237662603),
    @(21.53321237915288, 3.34136311495345, 4.266672007709356, 10.68491071602453, 72.12211101564753, 0, 7.344005520526261, 0, 11.30535622398139, 17.63169933834635, 10.97718810256054),
    @(21.53187863152725, 3.329147628686373, 4.256639557868513, 10.68740047215113, 72.0472195599708, 0, 7.344005520526261, 0, 11.30443295081869, 17.6323101391967, 10.98045342843031),
    @(21.53833778010304, 3.383752172526919, 4.301552240909416, 10.67628471199148, 72.38383971842831, 0, 7.344005520526261, 0, 11.3086159347095, 17.62992959314728, 10.96589149544294),
    @(21.56802690921943, 3.560796738178571, 4.448304377423186, 10.64048149964126, 73.50723891427073, 0, 7.344005520526261, 0, 11.3231616129018, 17.62837513878297, 10.91928087188522),
    @(21.59336469153099, 3.67171765872615, 4.541072885911563, 10.61822919393149, 74.23506134341241, 0, 7.344005520526261, 0, 11.33303275587654, 17.63216177656213, 10.89054215134617),
    @(21.57929985661834, 3.612955255812286, 4.491851159903608, 10.63000094307953, 73.84723863789617, 0, 7.344005520526261, 0, 11.32773107501583, 17.62970146513685, 10.90572290868881),
    @(21.53795298593617, 3.380777685609162, 4.299101318706056, 10.67688931736595, 72.36538087437886, 0, 7.344005520526261, 0, 11.30838437264616, 17.63003608363791, 10.96668245174571),
    @(21.5189378834119, 3.109816227802116, 4.078043950034412, 10.73241584290823, 70.74425514754175, 0, 7.344005520526261, 0, 11.28909136799074, 17.65100780803971, 11.03985065355964)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $r = $i + 2
    for ($j = 0; $j -lt $row.Length; $j++) {
        $c = $j + 2
        $ws.Cells.Item($r, $c).Value = $row[$j]
    }
}

Write-Output "updated"
